$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.211.45"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").Value = "2.367.35"
$ws.Range("E3").Value = "  +6.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.23"
$ws.Range("E5").Value = "  +3.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.98"
$ws.Range("E6").Value = "  -5.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +1.62%  "

# Row 8
$ws.Range("E8").Value = "  -0.30%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  +3.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.86"
$ws.Range("E10").Value = "  -5.60%  "

# Row 11
$ws.Range("E11").Value = "  +1.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.97"
$ws.Range("E12").Value = "  +0.70%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.07"
$ws.Range("E13").Value = "  +12.62%  "

# Row 14
$ws.Range("E14").Value = "  +1.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.56"
$ws.Range("E15").Value = "  +9.11%  "

# Row 16
$ws.Range("D16").Value = "2.725.09"
$ws.Range("E16").Value = "  +6.42%  "

# Row 17
$ws.Range("D17").Value = "2.379.08"
$ws.Range("E17").Value = "  +6.32%  "

# Row 18
$ws.Range("D18").Value = "43.198.50"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -0.01%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000109"
$ws.Range("E20").Value = "  +2.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.38"
$ws.Range("E21").Value = "  +2.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.43"
$ws.Range("E22").Value = "  -3.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.51"
$ws.Range("E23").Value = "  +8.80%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "252.20"
$ws.Range("E24").Value = "  +9.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.98"
$ws.Range("E25").Value = "  -5.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +2.09%  "

# Row 27
$ws.Range("E27").Value = "  +0.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.00"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("E29").Value = "  +1.13%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.72"
$ws.Range("E30").Value = "  +7.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.03"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32
$ws.Range("E32").Value = "  -2.30%  "

# Row 33
$ws.Range("E33").Value = "  +2.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  +1.01%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("E35").Value = "  +3.18%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("E36").Value = "  +1.24%  "

# Row 37
$ws.Range("E37").Value = "  +2.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.05"
$ws.Range("E38").Value = "  -5.81%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("E40").Value = "  +11.04%  "

# Row 41
$ws.Range("E41").Value = "  +15.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.35"
$ws.Range("E42").Value = "  +1.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.232"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.35"
$ws.Range("E45").Value = "  -6.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.68"
$ws.Range("E46").Value = "  +3.10%  "

# Row 47
$ws.Range("E47").Value = "  +9.75%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.39"
$ws.Range("E48").Value = "  +6.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").Value = "  -1.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0994"
$ws.Range("E50").Value = "  +0.97%  "

# Row 51
$ws.Range("D51").Value = "1.498.47"
$ws.Range("E51").Value = "  +4.45%  "
